# Generate Report for Handoff
# Updates the localization-status report with a newly generated handoff
# package: a new source-file GUID name and refreshed handoff timestamps /
# generated xliff file names for the zh-cn and de-de targets.

$wb = $excel.ActiveWorkbook

$oldFile = "ab27539e-faf2-4f64-bc16-850d27cc0899"
$newFile = "94c0e964-2323-4f37-ae93-3328c6e77f30"

$oldHash = "b0f4a7b8eff2c3dd8390f70980e2e370cc6e1518"
$newHash = "273af830355ba05e8f841c180ec6c5e3ce13173e"

# The hyperlinks all still point at the same commit-pinned source file in
# the repository (the relationship target itself does not change) - only
# the cell text / hyperlink display text is refreshed to the new name.
$mdUrl = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/346eabffddf7c2c54b914298889be696c622fcfc/e2e/$oldFile.md"

# NOTE: this PowerShell-style engine only binds function parameters
# positionally, so helper functions below are always invoked with
# positional arguments (no "-name value" syntax).
function Set-DisplayHyperlink {
    param($ws, [string]$cellRef, [string]$address, [string]$displayText)

    $rng = $ws.Range($cellRef)
    $rng.Value = $displayText

    # Re-create the hyperlink so its "display" text is refreshed while the
    # underlying link target/relationship is preserved.
    $rng.Hyperlinks.Delete()
    $ws.Hyperlinks.Add($rng, $address, [System.Reflection.Missing]::Value, [System.Reflection.Missing]::Value, $displayText) | Out-Null
}

# ---------------------------------------------------------------------
# Overview sheet
# ---------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")

$wsOverview.Range("A2").Value = "$newFile.md"

Set-DisplayHyperlink $wsOverview "B2" $mdUrl "e2e\$newFile.md"

$wsOverview.Range("G2").Value = "2016-08-22 20:57:37"

# ---------------------------------------------------------------------
# zh-cn sheet
# ---------------------------------------------------------------------
$wsZhCn = $wb.Worksheets.Item("zh-cn")

Set-DisplayHyperlink $wsZhCn "A2" $mdUrl "$newFile.md"

$wsZhCn.Range("G2").Value = "$newFile.$newHash.zh-cn.xlf"
$wsZhCn.Range("H2").Value = "2016-08-22 20:57:32"

# ---------------------------------------------------------------------
# de-de sheet
# ---------------------------------------------------------------------
$wsDeDe = $wb.Worksheets.Item("de-de")

Set-DisplayHyperlink $wsDeDe "A2" $mdUrl "$newFile.md"

$wsDeDe.Range("G2").Value = "$newFile.$newHash.de-de.xlf"
$wsDeDe.Range("H2").Value = "2016-08-22 20:57:37"
